# Applies the cell updates described in the commit diff for cryptos.xlsx
# (symbol list refresh performed by the GitHub Actions scraper).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell address -> new literal text value.
# All of these cells store plain text in the workbook (coin names, URLs,
# price strings and percentage strings), so each one is written back as
# text (avoiding Excel's automatic number/percentage conversion).
$updates = @(
    @{ Cell = 'D2'; Value = '307.63' }
    @{ Cell = 'E2'; Value = '-0.23%' }
    @{ Cell = 'D3'; Value = '41.03' }
    @{ Cell = 'E3'; Value = '0.91%' }
    @{ Cell = 'D4'; Value = '5.232' }
    @{ Cell = 'E4'; Value = '2.01%' }
    @{ Cell = 'D5'; Value = '0.07660' }
    @{ Cell = 'E5'; Value = '0.59%' }
    @{ Cell = 'B6'; Value = 'FTXToken' }
    @{ Cell = 'C6'; Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt' }
    @{ Cell = 'D6'; Value = '1.642' }
    @{ Cell = 'E6'; Value = '1.22%' }
    @{ Cell = 'B7'; Value = 'MXToken' }
    @{ Cell = 'C7'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D7'; Value = '0.9147' }
    @{ Cell = 'E7'; Value = '1.44%' }
    @{ Cell = 'B8'; Value = 'BTSEToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse' }
    @{ Cell = 'D8'; Value = '2.437' }
    @{ Cell = 'E8'; Value = '-0.45%' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D9'; Value = '0.1244' }
    @{ Cell = 'E9'; Value = '13.03%' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D10'; Value = '0.1825' }
    @{ Cell = 'E10'; Value = '3.81%' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D11'; Value = '0.09173' }
    @{ Cell = 'E11'; Value = '0.17%' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D12'; Value = '0.04159' }
    @{ Cell = 'E12'; Value = '-0.57%' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D13'; Value = '0.1051' }
    @{ Cell = 'E13'; Value = '0.00%' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D14'; Value = '0.001262' }
    @{ Cell = 'E14'; Value = '0.21%' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D15'; Value = '0.005832' }
    @{ Cell = 'E15'; Value = '-0.01%' }
    @{ Cell = 'B16'; Value = 'UpBots' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt' }
    @{ Cell = 'D16'; Value = '0.007509' }
    @{ Cell = 'E16'; Value = '2,395.62%' }
    @{ Cell = 'B17'; Value = 'LEO' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D17'; Value = '3.344' }
    @{ Cell = 'E17'; Value = '-0.28%' }
    @{ Cell = 'B18'; Value = 'GateToken' }
    @{ Cell = 'C18'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D18'; Value = '4.312' }
    @{ Cell = 'E18'; Value = '1.28%' }
    @{ Cell = 'D20'; Value = '7.463' }
    @{ Cell = 'E20'; Value = '13.23%' }
    @{ Cell = 'D21'; Value = '0.1393' }
    @{ Cell = 'E21'; Value = '2.02%' }
    @{ Cell = 'D22'; Value = '0.2884' }
    @{ Cell = 'E22'; Value = '7.55%' }
    @{ Cell = 'D23'; Value = '0.04057' }
    @{ Cell = 'E23'; Value = '-0.20%' }
    @{ Cell = 'D24'; Value = '0.001264' }
    @{ Cell = 'E24'; Value = '3.34%' }
    @{ Cell = 'D25'; Value = '0.004278' }
    @{ Cell = 'E25'; Value = '4.68%' }
    @{ Cell = 'D26'; Value = '0.0001273' }
    @{ Cell = 'E26'; Value = '-2.11%' }
    @{ Cell = 'D38'; Value = '0.02484' }
    @{ Cell = 'E38'; Value = '4.65%' }
    @{ Cell = 'D39'; Value = '0.05332' }
    @{ Cell = 'E39'; Value = '3.21%' }
    @{ Cell = 'D40'; Value = '0.007854' }
    @{ Cell = 'E40'; Value = '1.21%' }
    @{ Cell = 'E41'; Value = '1.02%' }
    @{ Cell = 'D42'; Value = '0.006589' }
    @{ Cell = 'E42'; Value = '-2.42%' }
    @{ Cell = 'E43'; Value = '-1.86%' }
    @{ Cell = 'D44'; Value = '0.007665' }
    @{ Cell = 'E44'; Value = '-12.52%' }
    @{ Cell = 'D45'; Value = '0.3347' }
    @{ Cell = 'E45'; Value = '0.35%' }
    @{ Cell = 'D46'; Value = '0.00006711' }
    @{ Cell = 'E46'; Value = '-4.42%' }
    @{ Cell = 'D47'; Value = '0.00000000752' }
    @{ Cell = 'E47'; Value = '0.20%' }
    @{ Cell = 'D48'; Value = '0.3705' }
    @{ Cell = 'E48'; Value = '1,075.09%' }
    @{ Cell = 'D49'; Value = '0.003107' }
    @{ Cell = 'E49'; Value = '-26.14%' }
    @{ Cell = 'D50'; Value = '0.00002105' }
    @{ Cell = 'E50'; Value = '0.20%' }
    @{ Cell = 'D51'; Value = '0.0002005' }
    @{ Cell = 'E51'; Value = '0.20%' }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    # Force text interpretation so numeric-looking / percent-looking
    # strings (e.g. '307.63', '-0.23%') are stored as text, matching
    # the original inline-string cell content instead of being
    # auto-converted into numbers by Excel.
    $cell.NumberFormat = '@'
    $cell.Value = $u.Value
    $cell.Style = 'Normal'
}

Write-Host "Updated $($updates.Count) cells"
